# Contra Costa County COVID testing data update.
# - Retroactive revisions to several historical days' "# of New Cases" (col C)
#   and "# of New Tests" (col F) ripple forward into the running totals
#   "# of Total Cases" (col B) and "# of Total Tests" (col E).
# - The trailing 7-day rolling figures (cols D and G) are re-published only
#   for the last handful of rows, matching the upstream source refresh.
# - A new day (2021-01-12, row 294) is appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> delta applied to that row's own "# of New Cases" (column C).
# Each such revision also shifts every later row's cumulative "# of Total
# Cases" (column B) by the same amount.
$cDeltas = @{
    163 = -1
    260 = 1
    282 = 1
    283 = 1
    288 = 1
    289 = -1
    290 = 1
}

# Row -> delta applied to that row's own "# of New Tests" (column F).
# Each such revision also shifts every later row's cumulative "# of Total
# Tests" (column E) by the same amount.
$fDeltas = @{
    34  = 1
    225 = -1
    238 = 1
    245 = 1
    255 = 2
    260 = 1
    261 = 1
    275 = 1
    281 = 1
    282 = 2
    287 = 7
    288 = 11
    289 = 2
    290 = 2
    291 = 1
    293 = 1
}

$runningB = 0
$runningE = 0

for ($r = 2; $r -le 293; $r++) {

    if ($cDeltas.ContainsKey($r)) {
        $cCur = $ws.Cells.Item($r, 3).Value2
        $ws.Cells.Item($r, 3).Value2 = $cCur + $cDeltas[$r]
        $runningB = $runningB + $cDeltas[$r]
    }

    if ($fDeltas.ContainsKey($r)) {
        $fCur = $ws.Cells.Item($r, 6).Value2
        $ws.Cells.Item($r, 6).Value2 = $fCur + $fDeltas[$r]
        $runningE = $runningE + $fDeltas[$r]
    }

    if ($runningB -ne 0) {
        $bCur = $ws.Cells.Item($r, 2).Value2
        $ws.Cells.Item($r, 2).Value2 = $bCur + $runningB
    }

    if ($runningE -ne 0) {
        $eCur = $ws.Cells.Item($r, 5).Value2
        $ws.Cells.Item($r, 5).Value2 = $eCur + $runningE
    }
}

# The upstream source also republished the trailing 7-day rolling average
# of new cases (D) and of % positive tests (G) for the most recent days.
$ws.Cells.Item(289, 4).Value2 = 655.285714285714
$ws.Cells.Item(289, 7).Value2 = 0.09140545603

$ws.Cells.Item(290, 4).Value2 = 633.857142857142
$ws.Cells.Item(290, 7).Value2 = 0.087423403542

$ws.Cells.Item(291, 4).Value2 = 597.285714285714
$ws.Cells.Item(291, 7).Value2 = 0.082795358232

$ws.Cells.Item(292, 4).Value2 = 585.714285714285
$ws.Cells.Item(292, 7).Value2 = 0.081494732657

$ws.Cells.Item(293, 4).Value2 = 579.857142857142
$ws.Cells.Item(293, 7).Value2 = 0.08075040783

# Append the new day at the bottom of the table.
$ws.Cells.Item(294, 1).Value2 = 44214
$ws.Cells.Item(294, 2).Value2 = 53678
$ws.Cells.Item(294, 3).Value2 = 463
$ws.Cells.Item(294, 4).Value2 = 525.285714285714
$ws.Cells.Item(294, 5).Value2 = 961714
$ws.Cells.Item(294, 6).Value2 = 6370
$ws.Cells.Item(294, 7).Value2 = 0.07847280022100001

Write-Output "Applied revisions and appended row 294."
